# Fixed point analysis methods paragraph - edit per commit
# "Minor additions to figures and edited methods paragraph".
#
# Two changes are applied:
#   1. The heading "Fixed point analysis" becomes bold.
#   2. The methods paragraph is substantially reworded/expanded.
#
# Both paragraphs are rewritten wholesale via Range.InsertXML (a genuine
# Word Range method) so the resulting run/text boundaries match the
# target OOXML exactly, rather than trying to replicate dozens of
# Find/Replace calls that each split runs in very specific places.

$d = $word.ActiveDocument

# --- Paragraph 1: bold the title run ("Fixed point analysis") ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Fixed point analysis</w:t></w:r></w:p>')

# --- Paragraph 2: replace with the revised/expanded methods text ---
$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">We </w:t></w:r><w:r><w:t>analyzed</w:t></w:r><w:r><w:t xml:space="preserve"> the</w:t></w:r><w:r><w:t xml:space="preserve"> model’s</w:t></w:r><w:r><w:t xml:space="preserve"> neural dynamics during the successful completion of a task. Following Driscoll et al., 2024, we used the Fixed Point Finder package [</w:t></w:r><w:r><w:t>https://github.com/mattgolub/fixed-point-finder</w:t></w:r><w:r><w:t>] to search for stationary points corresponding to correct model output.</w:t></w:r><w:r><w:t xml:space="preserve"> We stuck to the default arguments for Fixed Point Finder, following the example of the 3 bit flip-flop, including 5,000 maximum iterations and 1,000 initial states.</w:t></w:r><w:r><w:t xml:space="preserve"> We hypothesized that an arrangement of these features in state space </w:t></w:r><w:r><w:t>would</w:t></w:r><w:r><w:t xml:space="preserve"> characterize the trained models </w:t></w:r><w:r><w:t>and that this characterization would be</w:t></w:r><w:r><w:t xml:space="preserve"> invariant to the specific solutions (weight matrices) found during training. However, we did not find evidence of fixed points in our analysis. In five models </w:t></w:r><w:r><w:t>initialized with different random seeds</w:t></w:r><w:r><w:t xml:space="preserve">, we observed that the minimum q values were not sufficiently small </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t>&lt;</w:t></w:r><w:r><w:t xml:space="preserve">&lt;1/T^2 where T is the timescale of our task) </w:t></w:r><w:r><w:t xml:space="preserve">to be considered fixed points. To confirm that the state vector had velocity while the model output remained </w:t></w:r><w:r><w:t xml:space="preserve">sufficiently </w:t></w:r><w:r><w:t xml:space="preserve">constant </w:t></w:r><w:r><w:t xml:space="preserve">to meet </w:t></w:r><w:r><w:t xml:space="preserve">our </w:t></w:r><w:r><w:t xml:space="preserve">performance </w:t></w:r><w:r><w:t xml:space="preserve">criteria, we calculated the magnitude of the velocity vector over time from our hidden states. One can see </w:t></w:r><w:r><w:t xml:space="preserve">an appreciable velocity on the scale of sqrt(q). Cosine similarity between successive velocity vectors revealed ballistic motion during the beginning of the task followed by diffusive motion dominated by the intrinsic noise we include in our equations of motion. </w:t></w:r><w:r><w:t>W</w:t></w:r><w:r><w:t>hile the hidden state vector does not</w:t></w:r><w:r><w:t>, strictly speaking,</w:t></w:r><w:r><w:t xml:space="preserve"> settle into a stable fixed point as we hypothesized, it does fluctuate around a constant position. Given that the hidden state dimension is much greater than the output dimension, the output matrix is under</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">determined and multiple hidden states map to </w:t></w:r><w:r><w:t xml:space="preserve">similar enough </w:t></w:r><w:r><w:t>output state</w:t></w:r><w:r><w:t>s to meet our performance criteria</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>')
